# Scheduled data-refresh: push latest market-board derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns, H:N) into the
# per-job Fenrir_Profits sheets (ALC, ARM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 19095.5
$ws.Range("I28").Value = 430
$ws.Range("J28").Value = 32428
$ws.Range("K28").Value = 430
$ws.Range("L28").Value = 32428
$ws.Range("M28").Value = 55
$ws.Range("N28").Value = -33398
# Row 92
$ws.Range("H92").Value = 1370.762
$ws.Range("I92").Value = 777.3333
$ws.Range("J92").Value = 1815.8334
$ws.Range("K92").Value = 777.3333
$ws.Range("L92").Value = 1815.8334
$ws.Range("M92").Value = 470.6667
$ws.Range("N92").Value = -4311.8334
# Row 98
$ws.Range("H98").Value = 53423816
$ws.Range("I98").Value = 77166800
$ws.Range("J98").Value = 2100
$ws.Range("K98").Value = 77166800
$ws.Range("L98").Value = 2100
$ws.Range("M98").Value = -77165302
$ws.Range("N98").Value = -5096
# Row 122
$ws.Range("H122").Value = 53423816
$ws.Range("I122").Value = 77166800
$ws.Range("J122").Value = 2100
$ws.Range("K122").Value = 231500400
$ws.Range("L122").Value = 6300
$ws.Range("M122").Value = -231497950
$ws.Range("N122").Value = -11200
# Row 135
$ws.Range("H135").Value = 3729.0637
$ws.Range("I135").Value = 3914.8647
$ws.Range("J135").Value = 3041.6
$ws.Range("K135").Value = 35233.7823
$ws.Range("L135").Value = 27374.4
$ws.Range("M135").Value = -32698.7823
$ws.Range("N135").Value = -32444.4
# Row 137
$ws.Range("H137").Value = 21082296
$ws.Range("I137").Value = 331447.22
$ws.Range("K137").Value = 994341.6599999999
$ws.Range("M137").Value = -991791.6599999999
# Row 138
$ws.Range("H138").Value = 1773.6154
$ws.Range("I138").Value = 1292.0541
$ws.Range("J138").Value = 2409.9644
$ws.Range("K138").Value = 3876.1623
$ws.Range("L138").Value = 7229.8932
$ws.Range("M138").Value = 1263.8377
$ws.Range("N138").Value = -17509.8932

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4411.0625
$ws.Range("I32").Value = 3720.7642
$ws.Range("J32").Value = 13187.714
$ws.Range("K32").Value = 3720.7642
$ws.Range("L32").Value = 13187.714
$ws.Range("M32").Value = -3433.7642
$ws.Range("N32").Value = -13761.714
# Row 45
$ws.Range("H45").Value = 866.6
$ws.Range("I45").Value = 865.8182
$ws.Range("J45").Value = 868.75
$ws.Range("K45").Value = 865.8182
$ws.Range("L45").Value = 868.75
$ws.Range("M45").Value = -488.8182
$ws.Range("N45").Value = -1622.75
# Row 61
$ws.Range("H61").Value = 9527181
$ws.Range("I61").Value = 11908514
$ws.Range("J61").Value = 1852
$ws.Range("K61").Value = 11908514
$ws.Range("L61").Value = 1852
$ws.Range("M61").Value = -11908302
$ws.Range("N61").Value = -2276
# Row 74
$ws.Range("H74").Value = 585.561
$ws.Range("I74").Value = 381.74194
$ws.Range("J74").Value = 1217.4
$ws.Range("K74").Value = 381.74194
$ws.Range("L74").Value = 1217.4
$ws.Range("M74").Value = 492.25806
$ws.Range("N74").Value = -2965.4
# Row 77
$ws.Range("H77").Value = 585.561
$ws.Range("I77").Value = 381.74194
$ws.Range("J77").Value = 1217.4
$ws.Range("K77").Value = 1908.7097
$ws.Range("L77").Value = 6087
$ws.Range("M77").Value = 2459.2903
$ws.Range("N77").Value = -14823
# Row 136
$ws.Range("H136").Value = 9527181
$ws.Range("I136").Value = 11908514
$ws.Range("J136").Value = 1852
$ws.Range("K136").Value = 35725542
$ws.Range("L136").Value = 5556
$ws.Range("M136").Value = -35722992
$ws.Range("N136").Value = -10656

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 19519.254
$ws.Range("I58").Value = 21563.938
$ws.Range("J58").Value = 10597
$ws.Range("K58").Value = 21563.938
$ws.Range("L58").Value = 10597
$ws.Range("M58").Value = -21360.938
$ws.Range("N58").Value = -11003
# Row 134
$ws.Range("H134").Value = 16602683
$ws.Range("I134").Value = 19231622
$ws.Range("J134").Value = 5210618
$ws.Range("K134").Value = 57694866
$ws.Range("L134").Value = 15631854
$ws.Range("M134").Value = -57692331
$ws.Range("N134").Value = -15636924
# Row 136
$ws.Range("H136").Value = 19519.254
$ws.Range("I136").Value = 21563.938
$ws.Range("J136").Value = 10597
$ws.Range("K136").Value = 64691.814
$ws.Range("L136").Value = 31791
$ws.Range("M136").Value = -62141.814
$ws.Range("N136").Value = -36891

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 50001028
$ws.Range("I122").Value = 66667490
$ws.Range("J122").Value = 1642.2
$ws.Range("K122").Value = 200002470
$ws.Range("L122").Value = 4926.6
$ws.Range("M122").Value = -200000020
$ws.Range("N122").Value = -9826.6

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1430.1875
$ws.Range("I61").Value = 987.3
$ws.Range("J61").Value = 2168.3333
$ws.Range("K61").Value = 987.3
$ws.Range("L61").Value = 2168.3333
$ws.Range("M61").Value = -785.3
$ws.Range("N61").Value = -2572.3333
# Row 88
$ws.Range("H88").Value = 5000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -4572
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 5000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -3518
$ws.Range("N91").ClearContents()
# Row 113
$ws.Range("H113").Value = 1430.1875
$ws.Range("I113").Value = 987.3
$ws.Range("J113").Value = 2168.3333
$ws.Range("K113").Value = 987.3
$ws.Range("L113").Value = 2168.3333
$ws.Range("M113").Value = 1182.7
$ws.Range("N113").Value = -6508.3333
# Row 122
$ws.Range("H122").Value = 312501760
$ws.Range("I122").Value = 500000500
$ws.Range("J122").Value = 125003000
$ws.Range("K122").Value = 1500001500
$ws.Range("L122").Value = 375009000
$ws.Range("M122").Value = -1499999050
$ws.Range("N122").Value = -375013900
# Row 136
$ws.Range("H136").Value = 29416772
$ws.Range("I136").Value = 43484384
$ws.Range("J136").Value = 2673.0908
$ws.Range("K136").Value = 130453152
$ws.Range("L136").Value = 8019.2724
$ws.Range("M136").Value = -130450602
$ws.Range("N136").Value = -13119.2724

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 27793956
$ws.Range("I136").Value = 23828364
$ws.Range("J136").Value = 50001276
$ws.Range("K136").Value = 71485092
$ws.Range("L136").Value = 150003828
$ws.Range("M136").Value = -71482542
$ws.Range("N136").Value = -150008928
# Row 138
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280
